$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '69.876.57'
Set-TextValue "E2" '  -0.73%  '
Set-TextValue "D3" '3.583.60'
Set-TextValue "E3" '  -0.68%  '
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  -0.03%  '
Set-TextValue "D5" '579.41'
Set-TextValue "E5" '  -1.65%  '
Set-TextValue "D6" '190.60'
Set-TextValue "E6" '  +0.15%  '
Set-TextValue "D7" '0.632'
Set-TextValue "E7" '  -2.14%  '
Set-TextValue "D8" '3.580.16'
Set-TextValue "E8" '  -0.54%  '
Set-TextValue "D9" '1.00'
Set-TextValue "E9" '  +0.06%  '
Set-TextValue "D10" '0.180'
Set-TextValue "E10" '  +1.27%  '
Set-TextValue "D11" '0.663'
Set-TextValue "E11" '  +0.32%  '
Set-TextValue "D12" '55.70'
Set-TextValue "E12" '  -4.08%  '
Set-TextValue "D13" '0.0000305'
Set-TextValue "E13" '  +4.93%  '
Set-TextValue "D14" '9.63'
Set-TextValue "E14" '  -1.58%  '
Set-TextValue "D15" '4.155.70'
Set-TextValue "E15" '  -0.61%  '
Set-TextValue "D16" '19.84'
Set-TextValue "E16" '  +2.42%  '
Set-TextValue "D17" '3.579.17'
Set-TextValue "E17" '  -0.75%  '
Set-TextValue "D18" '69.868.55'
Set-TextValue "E18" '  -0.53%  '
Set-TextValue "D19" '12.65'
Set-TextValue "E19" '  +1.42%  '
Set-TextValue "D20" '0.121'
Set-TextValue "E20" '  +0.21%  '
Set-TextValue "D21" '1.04'
Set-TextValue "E21" '  -0.90%  '
Set-TextValue "D22" '473.94'
Set-TextValue "E22" '  -4.09%  '
Set-TextValue "D23" '19.48'
Set-TextValue "E23" '  +12.85%  '
Set-TextValue "D24" '5.02'
Set-TextValue "E24" '  -6.45%  '
Set-TextValue "D25" '95.74'
Set-TextValue "E25" '  +5.44%  '
Set-TextValue "D26" '4.37'
Set-TextValue "E26" '  -1.96%  '
Set-TextValue "D27" '3.00'
Set-TextValue "E27" '  -3.63%  '
Set-TextValue "D28" '11.02'
Set-TextValue "E28" '  -0.57%  '
Set-TextValue "D29" '9.28'
Set-TextValue "E29" '  -1.88%  '
Set-TextValue "D30" '32.19'
Set-TextValue "E30" '  -0.49%  '
Set-TextValue "D31" '7.63'
Set-TextValue "E31" '  +0.74%  '
Set-TextValue "D32" '12.20'
Set-TextValue "E32" '  -0.35%  '
Set-TextValue "D33" '0.119'
Set-TextValue "E33" '  +1.03%  '
Set-TextValue "D34" '66.32'
Set-TextValue "E34" '  +1.70%  '
Set-TextValue "D35" '590.31'
Set-TextValue "E35" '  -4.68%  '
Set-TextValue "D36" '38.96'
Set-TextValue "E36" '  +2.25%  '
Set-TextValue "D37" '1.00'
Set-TextValue "E37" '  +0.11%  '
Set-TextValue "D38" '0.0₃0799'
Set-TextValue "E38" '  -2.30%  '
Set-TextValue "D39" '0.395'
Set-TextValue "E39" '  -2.49%  '
Set-TextValue "D40" '3.18'
Set-TextValue "E40" '  +17.13%  '
Set-TextValue "D41" '3.46'
Set-TextValue "E41" '  -4.45%  '
Set-TextValue "D42" '0.137'
Set-TextValue "E42" '  -6.08%  '
Set-TextValue "D43" '3.221.33'
Set-TextValue "E43" '  -2.49%  '
Set-TextValue "D44" '2.85'
Set-TextValue "E44" '  +6.44%  '
Set-TextValue "D45" '3.07'
Set-TextValue "E45" '  -0.74%  '
Set-TextValue "D46" '0.0443'
Set-TextValue "E46" '  -0.57%  '
Set-TextValue "D47" '3.36'
Set-TextValue "E47" '  +1.82%  '
Set-TextValue "D48" '9.43'
Set-TextValue "E48" '  +3.17%  '
Set-TextValue "D49" '0.138'
Set-TextValue "E49" '  -0.03%  '
Set-TextValue "D50" '0.998'
Set-TextValue "E50" '  -0.17%  '
Set-TextValue "D51" '3.14'
Set-TextValue "E51" '  -5.42%  '
